$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValuationTable")
$bi = $wb.Worksheets.Item("BuildingInfoTable")

$xlPasteFormats = -4122

function CopyFormat($srcRange, $dstRange) {
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial($xlPasteFormats) | Out-Null
}

# --- New column K formatting (copy from column J in same row), for rows without new values ---
CopyFormat $ws.Range("J1") $ws.Range("K1")
CopyFormat $ws.Range("J3") $ws.Range("K3")
CopyFormat $ws.Range("J4") $ws.Range("K4")
CopyFormat $ws.Range("J7") $ws.Range("K7")
CopyFormat $ws.Range("J8") $ws.Range("K8")
CopyFormat $ws.Range("J9") $ws.Range("K9")
CopyFormat $ws.Range("J10") $ws.Range("K10")
CopyFormat $ws.Range("J11") $ws.Range("K11")

Write-Host "done stage1"

# --- Row 2 header: insert "agriculturalMarketValue" at F2, shifting G2:J2 text right
#     into G2:K2. Source cells A2/C2 = style 3 ; B2/D2/E2/F2/G2 = style 4 (untouched by
#     this edit, safe copy sources) ---
CopyFormat $ws.Range("C2") $ws.Range("K2")
CopyFormat $ws.Range("G2") $ws.Range("H2")
CopyFormat $ws.Range("C2") $ws.Range("I2")
CopyFormat $ws.Range("G2") $ws.Range("J2")

$ws.Range("F2").Value2 = "agriculturalMarketValue"
$ws.Range("G2").Value2 = "justMarketValue"
$ws.Range("H2").Value2 = "assessedValue"
$ws.Range("I2").Value2 = "exemptValue"
$ws.Range("J2").Value2 = "taxableValue"
$ws.Range("K2").Value2 = "maximumSaveOurHomesPortability"

Write-Host "done stage2"

# --- Row 5: populate data-source row, styles borrowed from BuildingInfoTable row 5 ---
# (B5/D5/E5/F5/G5/I5 need the "1,3" string styling; C5/H5/J5 already default to the
#  correct numeric style; K5 is a brand new cell and needs the border format copied in)
CopyFormat $bi.Range("B5") $ws.Range("B5")
CopyFormat $bi.Range("C5") $ws.Range("D5")
CopyFormat $bi.Range("C5") $ws.Range("E5")
CopyFormat $bi.Range("C5") $ws.Range("F5")
CopyFormat $bi.Range("C5") $ws.Range("G5")
CopyFormat $bi.Range("C5") $ws.Range("I5")
CopyFormat $ws.Range("C4") $ws.Range("K5")

$ws.Range("B5").Value2 = "1,3"
$ws.Range("D5").Value2 = "1,3"
$ws.Range("E5").Value2 = "1,3"
$ws.Range("F5").Value2 = "1,3"
$ws.Range("G5").Value2 = "1,3"
$ws.Range("I5").Value2 = "1,3"
$ws.Range("C5").Value2 = 3
$ws.Range("H5").Value2 = 3
$ws.Range("J5").Value2 = 3
$ws.Range("K5").Value2 = 0

Write-Host "done stage3"

# --- Row 6: new values, add agriculturalMarketValue column, shift remaining ---
# Styles before: A8 B30 C15 D31 E31 F31 G31 H32 I31 J15
# Styles after:  A8 B30 C15 D31 E31 F15 G31 H31 I32 J31 K15
# I needs H's ORIGINAL style (32) - copy first, before H changes.
CopyFormat $ws.Range("H6") $ws.Range("I6")
# H, J need style 31 (borrow from D/E/G which remain 31)
CopyFormat $ws.Range("D6") $ws.Range("H6")
CopyFormat $ws.Range("E6") $ws.Range("J6")
# F, K need style 15 (borrow from C which remains 15)
CopyFormat $ws.Range("C6") $ws.Range("F6")
CopyFormat $ws.Range("C6") $ws.Range("K6")

$ws.Range("B6").Value2 = "justvalhomestead"
$ws.Range("D6").Value2 = "lnd_val"
$ws.Range("E6").Value2 = "asvalagval"
$ws.Range("F6").Value2 = "DATA NOT AVAILABLE"
$ws.Range("G6").Value2 = "parval"
$ws.Range("H6").Value2 = "assedvalschool"
$ws.Range("I6").Value2 = "assedvalnonschool"
$ws.Range("J6").Value2 = "taxvalschool"
$ws.Range("K6").Value2 = "DATA NOT AVAILABLE"

Write-Host "done stage4"

# --- Merge the new A1:K1 header cell (was A1:J1) ---
$ws.Range("A1:J1").UnMerge() | Out-Null
$ws.Range("A1:K1").Merge() | Out-Null

Write-Host "done stage5"

# --- Column widths: insert a new "agriculturalMarketValue" column (F) at width 20.5,
#     shifting the old "exemptValue" column (previously H, width 17.6406) to I ---
# Note: ColumnWidth values get snapped to the engine's internal character-width grid,
# so these inputs are chosen to land as close as possible to the target stored widths.
$ws.Columns.Item(6).ColumnWidth = 19.8    # -> stored width ~20.57 (closest to target 20.5)
$ws.Columns.Item(8).ColumnWidth = 15.57   # -> reset old H column back to the default width
$ws.Columns.Item(9).ColumnWidth = 16.86   # -> stored width ~17.57 (closest to target 17.6406)

Write-Host "done stage6"
